# C5-PowerPoint.pptx edit:
#  1) Slide 6's table switches from table style {FFA53DC7-...} to {71FE0D19-...}
#  2) The presentation's theme (ppt/theme/theme2.xml, used by the slide master /
#     every slide) switches its colour palette from "Integral" to the stock
#     "Office Theme" palette (font scheme / format scheme are identical between
#     the two themes already, only the 12 colour-scheme entries differ).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 --------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{71FE0D19-2B52-43A7-A2E3-551EBF480797}")

# --- 2) Theme colour scheme: Integral -> Office Theme ----------------------
function Set-ThemeRGB($scheme, $index, $r, $g, $b) {
    $scheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$themeColors = $slide.ThemeColorScheme
Set-ThemeRGB $themeColors 1  0   0   0    # dk1
Set-ThemeRGB $themeColors 2  255 255 255  # lt1
Set-ThemeRGB $themeColors 3  68  84  106  # dk2      -> 44546A
Set-ThemeRGB $themeColors 4  231 230 230  # lt2      -> E7E6E6
Set-ThemeRGB $themeColors 5  91  155 213  # accent1  -> 5B9BD5
Set-ThemeRGB $themeColors 6  237 125 49   # accent2  -> ED7D31
Set-ThemeRGB $themeColors 7  165 165 165  # accent3  -> A5A5A5
Set-ThemeRGB $themeColors 8  255 192 0    # accent4  -> FFC000
Set-ThemeRGB $themeColors 9  68  114 196  # accent5  -> 4472C4
Set-ThemeRGB $themeColors 10 112 173 71   # accent6  -> 70AD47
Set-ThemeRGB $themeColors 11 5   99  193  # hlink    -> 0563C1
Set-ThemeRGB $themeColors 12 149 79  114  # folHlink -> 954F72
